$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (co2)
$ws.Range("C2").Value = 552.5041813615475
$ws.Range("D2").Value = 133.2343220950599
$ws.Range("F2").Value = 448
$ws.Range("G2").Value = 517
$ws.Range("H2").Value = 619

# Row 3 (humidity)
$ws.Range("C3").Value = 41.64199416212463
$ws.Range("D3").Value = 4.854508218915211
$ws.Range("F3").Value = 38.55
$ws.Range("G3").Value = 40.66
$ws.Range("H3").Value = 44.65

# Row 4 (pm25)
$ws.Range("C4").Value = 1.350393341883082
$ws.Range("D4").Value = 2.332899451123587
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.98
$ws.Range("H4").Value = 1.73

# Row 5 (pressure)
$ws.Range("C5").Value = 321.9401196296889
$ws.Range("D5").Value = 10.27790085095419
$ws.Range("F5").Value = 315.36
$ws.Range("G5").Value = 323.87
$ws.Range("H5").Value = 330.37

# Row 6 (temperature)
$ws.Range("C6").Value = 21.52157452976382
$ws.Range("D6").Value = 2.003594379020124
$ws.Range("F6").Value = 20.05
$ws.Range("G6").Value = 21.36
$ws.Range("H6").Value = 22.47

# Row 7 (rssi)
$ws.Range("C7").Value = -76.19348081650101
$ws.Range("D7").Value = 23.09578817080914

# Row 8 (snr)
$ws.Range("C8").Value = 7.686607516620905
$ws.Range("D8").Value = 6.89608249706923

# Row 9 (SF)
$ws.Range("C9").Value = 9.322246713690285
$ws.Range("D9").Value = 1.688607784523832

# Row 10 (frequency)
$ws.Range("C10").Value = 867.8303182109649
$ws.Range("D10").Value = 0.4611200429001951

# Row 11 (toa)
$ws.Range("C11").Value = 0.5569688428342418
$ws.Range("D11").Value = 0.5908372583519632

# Row 12 (distance)
$ws.Range("C12").Value = 22.70849898471732
$ws.Range("D12").Value = 12.28637657221718

# Row 13 (c_walls)
$ws.Range("C13").Value = 0.6722988137223469
$ws.Range("D13").Value = 0.7488477031850275

# Row 14 (w_walls)
$ws.Range("C14").Value = 1.827625040076948
$ws.Range("D14").Value = 1.667054721662187

# Row 15 (exp_pl)
$ws.Range("C15").Value = 93.59348081650087
$ws.Range("D15").Value = 23.09578817080798

# Row 16 (n_power)
$ws.Range("C16").Value = -85.45703572767513
$ws.Range("D16").Value = 20.85876560659655
$ws.Range("F16").Value = -101.8707776445072
$ws.Range("G16").Value = -85.0778545523916

# Row 17 (esp)
$ws.Range("C17").Value = -77.77042821105421
$ws.Range("D17").Value = 25.46089467460196
$ws.Range("F17").Value = -92.53779541063678
$ws.Range("G17").Value = -74.46183611348224
$ws.Range("H17").Value = -54.79009749652566
